# Atualização de bases das ligas, do dia: 09-03-2024 às 13:07
#
# 1) Rows 2 and 3 (match ids 6720873 / 6720843) swap places in the sheet.
# 2) Rows 144 and 145 (match ids 7493311 / 7493312) swap places in the sheet.
# 3) Rows 192 and 193 get their results (FTHG/FTAG/FTR) filled in, plus the
#    odds columns that move once a result is known.
# 4) Rows 194 and 195 get updated closing-line odds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap full data rows 2 <-> 3 (columns B..AC; A/id column stays put) ---
$row2 = $ws.Range("B2:AC2").Value()
$row3 = $ws.Range("B3:AC3").Value()
$ws.Range("B2:AC2").Value = $row3
$ws.Range("B3:AC3").Value = $row2

# --- 2) Swap full data rows 144 <-> 145 (columns B..AC; A/id column stays put) ---
$row144 = $ws.Range("B144:AC144").Value()
$row145 = $ws.Range("B145:AC145").Value()
$ws.Range("B144:AC144").Value = $row145
$ws.Range("B145:AC145").Value = $row144

# --- 3) Row 192: match finished 1-0 (home win) ---
$ws.Range("H192").Value = 1
$ws.Range("I192").Value = 0
$ws.Range("J192").Value = "H"
$ws.Range("N192").Value = 2.25
$ws.Range("P192").Value = 3.2
$ws.Range("R192").Value = 2.025
$ws.Range("S192").Value = 1.775
$ws.Range("W192").Value = 1.25
$ws.Range("X192").Value = -1
$ws.Range("Y192").Value = -1
$ws.Range("Z192").Value = 1.025
$ws.Range("AA192").Value = -1
$ws.Range("AB192").Value = -1
$ws.Range("AC192").Value = 0.8500000000000001

# --- Row 193: match finished 1-1 (draw) ---
$ws.Range("H193").Value = 1
$ws.Range("I193").Value = 1
$ws.Range("J193").Value = "D"
$ws.Range("N193").Value = 2.1
$ws.Range("O193").Value = 3.1
$ws.Range("P193").Value = 3.4
$ws.Range("R193").Value = 1.8
$ws.Range("S193").Value = 2
$ws.Range("W193").Value = -1
$ws.Range("X193").Value = 2.1
$ws.Range("Y193").Value = -1
$ws.Range("Z193").Value = -0.5
$ws.Range("AA193").Value = 0.5
$ws.Range("AB193").Value = -0.5
$ws.Range("AC193").Value = 0.3875

# --- 4) Row 194: updated closing odds (match still upcoming) ---
$ws.Range("N194").Value = 1.3
$ws.Range("O194").Value = 5
$ws.Range("P194").Value = 8.5
$ws.Range("Q194").Value = -1.5
$ws.Range("R194").Value = 1.9
$ws.Range("S194").Value = 1.9
$ws.Range("U194").Value = 1.8
$ws.Range("V194").Value = 2

# --- Row 195: updated closing odds (match still upcoming) ---
$ws.Range("N195").Value = 6.5
$ws.Range("O195").Value = 3.6
$ws.Range("P195").Value = 1.5
$ws.Range("R195").Value = 1.925
$ws.Range("S195").Value = 1.875
$ws.Range("U195").Value = 1.85
$ws.Range("V195").Value = 1.95
